# Generate Report for Handoff
#
# The "b.md" file has completed its handback cycle: its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" on every sheet,
# a new (non-latest) handback xliff file is recorded for both locales, the
# handback timestamps are refreshed, and an "out of date" warning is written
# into the Error Detail column. Column P (Error Detail) is also widened.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3a95437b7ed69af95b67a93e7bad94afbf1bd960/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b0f5828955dcaa1d12baac67e10bfa1405ad6d52/e2e/b.md."

# --- Overview sheet: row 3 is b.md ---
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = "2016-09-02 14:46:56"

# --- zh-cn sheet: row 3 is b.md ---
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-02 14:46:51"
$wsZhCn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 is b.md ---
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-02 14:46:56"
$wsDeDe.Range("P3").Value = $errorDetail

# --- Widen "Error Detail" column (P, column 16) on both locale sheets ---
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
